$d = $word.ActiveDocument

# The opening "S / ervice Interface / version:" heading is the sole
# paragraph using pStyle "Heading3" in this template; switch it to the
# (new) "No Spacing" style so the empty block under each domGroup collapses
# instead of reserving Heading-3 space.
$p1 = $d.Paragraphs.Item(1)
$p1.Style = "No Spacing"

# "No Spacing" did not previously exist in this style sheet, so assigning
# it minted a bare definition. Flesh it out to match Word's real built-in
# "No Spacing" style (uiPriority 1, Quick Style, Arial 10pt/no-proof,
# zero spacing-after / single line spacing, en-GB / de-DE language) which
# is exactly what Word itself writes the first time the style is used.
$s = $d.Styles.Item("No Spacing")
$s.Priority = 1
$s.QuickStyle = $true
$s.NoProofing = $true
$s.Font.Name = "Arial"
$s.Font.Size = 10
$s.Font.NameBi = "Arial"
$s.Font.SizeBi = 10
$s.Font.LanguageID = "en-GB"
$s.Font.LanguageIDFarEast = "de-DE"
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.LineSpacingRule = 0
